# The commit swaps the contents of ppt/theme/theme1.xml (the slide
# master's theme, originally the "Integral" theme) and
# ppt/theme/theme2.xml (the notes master's theme, originally the
# default "Office Theme"). Font scheme and format scheme (fills,
# lines, effects, background fills) are byte-identical between the
# two parts - only each theme's <a:clrScheme> (and the name
# attributes, which PowerPoint's object model does not expose as
# writable) differ. So the net visible effect is: the slide master's
# theme colors change from the "Integral" palette to the stock
# "Office" palette.
#
# Apply that by rewriting every entry of the (reachable) slide
# master's ThemeColorScheme to the 12 stock Office theme colors, in
# the fixed dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink order that
# ThemeColorScheme.Item(1..12) walks.
#
# ColorFormat.RGB uses the Windows COLORREF byte order (0x00BBGGRR),
# i.e. the bytes of the RRGGBB hex string reversed, so each literal
# below is the target srgbClr value byte-swapped.

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$scheme = $theme.ThemeColorScheme

$scheme.Item(1).RGB  = 0x000000   # dk1      -> srgbClr 000000
$scheme.Item(2).RGB  = 0xFFFFFF   # lt1      -> srgbClr FFFFFF
$scheme.Item(3).RGB  = 0x6A5444   # dk2      -> srgbClr 44546A
$scheme.Item(4).RGB  = 0xE6E6E7   # lt2      -> srgbClr E7E6E6
$scheme.Item(5).RGB  = 0xD59B5B   # accent1  -> srgbClr 5B9BD5
$scheme.Item(6).RGB  = 0x317DED   # accent2  -> srgbClr ED7D31
$scheme.Item(7).RGB  = 0xA5A5A5   # accent3  -> srgbClr A5A5A5
$scheme.Item(8).RGB  = 0x00C0FF   # accent4  -> srgbClr FFC000
$scheme.Item(9).RGB  = 0xC47244   # accent5  -> srgbClr 4472C4
$scheme.Item(10).RGB = 0x47AD70   # accent6  -> srgbClr 70AD47
$scheme.Item(11).RGB = 0xC16305   # hlink    -> srgbClr 0563C1
$scheme.Item(12).RGB = 0x724F95   # folHlink -> srgbClr 954F72
